$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = "stop"
$ws.Range("C47").Value = "stop"
$ws.Range("C55").Value = "stop"
$ws.Range("B62").Value = "stop"
$ws.Range("B69").Value = "stop"
$ws.Range("B71").Value = "stop"
$ws.Range("B74").Value = "stop"
$ws.Range("B75").Value = "stop"
$ws.Range("B76").Value = "stop"
$ws.Range("B77").Value = "stop"
$ws.Range("C80").Value = "other"
$ws.Range("C85").Value = "other"
$ws.Range("C87").Value = "other"
$ws.Range("C92").Value = "other"
$ws.Range("C94").Value = "other"
$ws.Range("C99").Value = "other"
$ws.Range("C110").Value = "other"
$ws.Range("C127").Value = "stop"
$ws.Range("C129").Value = "stop"
$ws.Range("C140").Value = "speedlimit"
$ws.Range("C145").Value = "other"
$ws.Range("C171").Value = "speedlimit"
$ws.Range("C190").Value = "other"
$ws.Range("C191").Value = "speedlimit"
$ws.Range("C196").Value = "other"
$ws.Range("C197").Value = "speedlimit"
$ws.Range("B198").Value = "stop"
$ws.Range("C198").Value = "stop"
$ws.Range("C203").Value = "stop"
$ws.Range("C204").Value = "stop"
$ws.Range("C205").Value = "stop"
$ws.Range("C206").Value = "stop"
$ws.Range("B209").Value = "stop"
$ws.Range("C209").Value = "stop"
$ws.Range("C210").Value = "speedlimit"
$ws.Range("B219").Value = "stop"
$ws.Range("B233").Value = "stop"
$ws.Range("C233").Value = "stop"
$ws.Range("B239").Value = "stop"
$ws.Range("B243").Value = "stop"
$ws.Range("C243").Value = "stop"
$ws.Range("B253").Value = "stop"
$ws.Range("C253").Value = "stop"
$ws.Range("C261").Value = "other"
$ws.Range("B264").Value = "stop"
$ws.Range("C264").Value = "stop"
$ws.Range("B272").Value = "stop"
$ws.Range("C272").Value = "stop"
$ws.Range("C276").Value = "speedlimit"
$ws.Range("B282").Value = "stop"
$ws.Range("C282").Value = "stop"
$ws.Range("C286").Value = "speedlimit"
$ws.Range("C287").Value = "speedlimit"
$ws.Range("B293").Value = "stop"
$ws.Range("C293").Value = "stop"
$ws.Range("B294").Value = "stop"
$ws.Range("C294").Value = "stop"
$ws.Range("B295").Value = "stop"
$ws.Range("C295").Value = "stop"
$ws.Range("B296").Value = "stop"
$ws.Range("C296").Value = "stop"
$ws.Range("B297").Value = "stop"
$ws.Range("C297").Value = "stop"
$ws.Range("B298").Value = "stop"
$ws.Range("C298").Value = "stop"
$ws.Range("B299").Value = "stop"
$ws.Range("C299").Value = "stop"
$ws.Range("B300").Value = "stop"
$ws.Range("B303").Value = "stop"
$ws.Range("C303").Value = "stop"
$ws.Range("C311").Value = "speedlimit"
$ws.Range("B313").Value = "stop"
$ws.Range("C313").Value = "stop"
$ws.Range("B316").Value = "stop"
$ws.Range("C323").Value = "stop"
$ws.Range("C325").Value = "speedlimit"
$ws.Range("C326").Value = "stop"
$ws.Range("B327").Value = "stop"
$ws.Range("C327").Value = "stop"
$ws.Range("B329").Value = "stop"
$ws.Range("B332").Value = "stop"
$ws.Range("B333").Value = "stop"
$ws.Range("C334").Value = "speedlimit"
$ws.Range("C335").Value = "speedlimit"
$ws.Range("B337").Value = "stop"
$ws.Range("C337").Value = "stop"
$ws.Range("B347").Value = "stop"
$ws.Range("C347").Value = "stop"
$ws.Range("C348").Value = "speedlimit"
$ws.Range("C356").Value = "speedlimit"
$ws.Range("B363").Value = "stop"
$ws.Range("C363").Value = "stop"
$ws.Range("C370").Value = "stop"
$ws.Range("B381").Value = "stop"
$ws.Range("C381").Value = "stop"
$ws.Range("C384").Value = "stop"
$ws.Range("B389").Value = "stop"
$ws.Range("B400").Value = "stop"
$ws.Range("C400").Value = "stop"
$ws.Range("B411").Value = "stop"
$ws.Range("C411").Value = "stop"
$ws.Range("C419").Value = "stop"
$ws.Range("C420").Value = "stop"
$ws.Range("B421").Value = "stop"
$ws.Range("C421").Value = "stop"
$ws.Range("C422").Value = "stop"
$ws.Range("C423").Value = "stop"
$ws.Range("C428").Value = "stop"
$ws.Range("B440").Value = "stop"
$ws.Range("C440").Value = "stop"
$ws.Range("B451").Value = "stop"
$ws.Range("C451").Value = "stop"
$ws.Range("B462").Value = "stop"
$ws.Range("C462").Value = "stop"
$ws.Range("B469").Value = "stop"
$ws.Range("C469").Value = "stop"
$ws.Range("C470").Value = "stop"
$ws.Range("C471").Value = "stop"
$ws.Range("C472").Value = "stop"
$ws.Range("C473").Value = "stop"
$ws.Range("C474").Value = "stop"
$ws.Range("B478").Value = "stop"
$ws.Range("C478").Value = "stop"
$ws.Range("C484").Value = "speedlimit"
$ws.Range("C485").Value = "speedlimit"
$ws.Range("B489").Value = "stop"
$ws.Range("C489").Value = "stop"
$ws.Range("C490").Value = "speedlimit"
$ws.Range("B500").Value = "stop"
$ws.Range("C500").Value = "stop"
$ws.Range("B521").Value = "stop"
$ws.Range("C521").Value = "stop"
$ws.Range("B524").Value = "stop"
$ws.Range("B525").Value = "stop"
$ws.Range("B526").Value = "stop"
$ws.Range("B527").Value = "stop"
$ws.Range("B528").Value = "stop"
$ws.Range("B529").Value = "stop"
$ws.Range("B530").Value = "stop"
$ws.Range("C530").Value = "stop"
$ws.Range("B531").Value = "stop"
$ws.Range("C531").Value = "stop"
$ws.Range("B532").Value = "stop"
$ws.Range("C532").Value = "stop"
$ws.Range("B533").Value = "stop"
$ws.Range("C533").Value = "stop"
$ws.Range("B534").Value = "stop"
$ws.Range("C534").Value = "stop"
$ws.Range("B535").Value = "stop"
$ws.Range("C535").Value = "stop"
$ws.Range("B536").Value = "stop"
$ws.Range("C536").Value = "stop"
